$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update odds on row 2 (ITALY - SERIE A, Genoa vs Fiorentina)
$ws.Range("O2").Value = 1.3
$ws.Range("P2").Value = 3.5
$ws.Range("Q2").Value = 2.04
$ws.Range("R2").Value = 1.86

# 2) Update odds on row 3 (ITALY - SERIE A, AS Roma vs Torino)
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 1.8

# 3) Remove the EGYPT - PREMIER LEAGUE fixture (old row 5); rows below shift up
$ws.Range("A5:BD5").EntireRow.Delete()

# 4) Tweak a handful of odds on the fixture that is now row 5
#    (SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE, Al Ittihad vs Al Ahli SC)
$ws.Range("H5").Value = 3.6
$ws.Range("M5").Value = 19
$ws.Range("N5").Value = 1.03
$ws.Range("AB5").Value = 21
$ws.Range("AM5").Value = 23
$ws.Range("AX5").Value = 15
$ws.Range("BB5").Value = 101
